$wb = $excel.ActiveWorkbook

# Rename the "CompanyRecordTypes" sheet to "RoundTripFields"
$ws = $wb.Worksheets.Item("CompanyRecordTypes")
$ws.Name = "RoundTripFields"

# Replace the sheet contents with the new Round Trip section field labels
$ws.Range("A1").Value = "Round Trip Section Fields"
$ws.Range("A2").Value = "Potential Round Trip"
$ws.Range("A3").Value = "Round Trip Engagement"
$ws.Range("A4").Value = "Round Trip Comment"
$ws.Range("A5").Value = "Potential Round Trip Last Modified Date"

# Widen column A to fit the new (longer) content
$ws.Columns.Item(1).ColumnWidth = 33.276041666666664

# Update the selected cell on this sheet
$ws.Range("A11").Select() | Out-Null

# Adjust the workbook window size/position as recorded in the saved view
$win = $wb.Windows.Item(1)
$win.Left = 3624
$win.Top = 1500
$win.Width = 19416
$win.Height = 10740
